# Add the new "multi line cells" worksheet as the last tab (after "strings"),
# matching sheetId=3 / rId3 ordering from the target workbook.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "multi line cells"

# --- Cell values -----------------------------------------------------
# Column B first, then C, then A/B4 so the shared-string table gets built
# in the same order as the authored workbook.
$ws.Range("B1").Value = "line 1 line 2 line 3"
$ws.Range("B2").Value = "line 1 `nline 2 `nline 3"
$ws.Range("B3").Value = "line 1`nline 2`nline 3"
$ws.Range("C4").Value = "wrap text"
$ws.Range("C2").Value = "alt+enter`nw/spaces"
$ws.Range("C3").Value = "alt+enter`nno spaces"

$ws.Range("A1").Value = "pass"
$ws.Range("A2").Value = "pass"
$ws.Range("A3").Value = "fail"
$ws.Range("A4").Value = "pass"
$ws.Range("B4").Value = "line 1 line 2 line 3"

# --- Formatting --------------------------------------------------------
$ws.Range("B2:B4").WrapText = $true
$ws.Range("C2:C3").WrapText = $true
$ws.Range("C2").VerticalAlignment = -4108   # xlVAlignCenter
$ws.Range("F3").WrapText = $true            # style-only cell, no value

$ws.Rows.Item(2).RowHeight = 43.2
$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 28.8

$ws.Columns.Item(2).ColumnWidth = 10.65

# --- View state ----------------------------------------------------
$ws.Range("D8").Select() | Out-Null

Write-Output "done"
